$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '38.268.58'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.49%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.121.10'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.80%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.94%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.626'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.98%  '

$ws.Range("E7").Value = '  +2.65%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("E9").Value = '  +3.66%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0782'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.13%  '

$ws.Range("E11").Value = '  +1.80%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.435.93'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.74%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.59'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.80%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.64'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.75%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.792'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.34%  '

$ws.Range("E16").Value = '  +2.93%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.121.49'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.67%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '38.188.25'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.36%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.17'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.63'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.96%  '

$ws.Range("E21").Value = '  +3.27%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.31'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.30%  '

$ws.Range("E23").Value = '  -0.01%  '

$ws.Range("E24").Value = '  -0.32%  '

$ws.Range("E25").Value = '  +2.57%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.84'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.00%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.140'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +12.80%  '

$ws.Range("E28").Value = '  +3.62%  '

$ws.Range("E29").Value = '  -0.93%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.58'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.46%  '

$ws.Range("E31").Value = '  +1.53%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.64'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.62%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.64'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.01%  '

$ws.Range("E34").Value = '  +2.96%  '

$ws.Range("E35").Value = '  +2.14%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.50'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.85%  '

$ws.Range("E37").Value = '  +4.68%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.01%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.50'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.25%  '

$ws.Range("E40").Value = '  +9.26%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.96'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.51%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '97.42'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.87%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0215'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.40%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.462.12'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.92%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.16'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.92%  '

$ws.Range("E46").Value = '  +6.46%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.14'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.29%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.79'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.28%  '

$ws.Range("E49").Value = '  +3.91%  '

$ws.Range("E50").Value = '  +3.46%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.320.65'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.74%  '
